$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-preserving number format on D-column (price) cells so
# numeric-looking strings like "1.00" or "595.67" are not coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "68.259.67"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "2.643.59"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "595.67"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "159.84"
$ws.Range("E6").Value = "  +3.50%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").Value = "27.90"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "3.126.15"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").Value = "68.098.29"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "2.645.02"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").Value = "360.04"
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").Value = "7.32"
$ws.Range("E20").Value = "  -2.30%  "
$ws.Range("D21").Value = "4.40"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "4.74"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "75.07"
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "9.78"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").Value = "2.812.80"
$ws.Range("E27").Value = "  +1.36%  "
$ws.Range("E28").Value = "  -3.37%  "
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").Value = "558.30"
$ws.Range("E29").Value = "  -2.84%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "0.918"
$ws.Range("E30").Value = "  -8.05%  "
$ws.Range("D31").Value = "7.97"
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -3.23%  "
$ws.Range("E36").Value = "  -2.38%  "
$ws.Range("D38").Value = "158.70"
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("E40").Value = "  -2.68%  "
$ws.Range("D41").Value = "5.32"
$ws.Range("E41").Value = "  -1.93%  "
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("E43").Value = "  -6.70%  "
$ws.Range("D45").Value = "156.82"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").Value = "21.99"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("E48").Value = "  -2.70%  "
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "0.565"
$ws.Range("E51").Value = "  -0.21%  "

# Restore default (Normal) style on D-column cells so no stray style index is left behind
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D51").Style = "Normal"

